$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sample names to reflect the leakage-test naming convention ---
# Mix1 / Mix2 / 3N2O / 10N2O / 3KCO2 (rows 2-16, first 2014.06.05 batch) -> *Leak
$ws.Range("E2").Value  = "Mix1Leak"
$ws.Range("E3").Value  = "Mix1Leak"
$ws.Range("E4").Value  = "Mix1Leak"
$ws.Range("E5").Value  = "Mix2Leak"
$ws.Range("E6").Value  = "Mix2Leak"
$ws.Range("E7").Value  = "Mix2Leak"
$ws.Range("E8").Value  = "3N2OLeak"
$ws.Range("E9").Value  = "3N2OLeak"
$ws.Range("E10").Value = "3N2OLeak"
$ws.Range("E11").Value = "10N2OLeak"
$ws.Range("E12").Value = "10N2OLeak"
$ws.Range("E13").Value = "10N2OLeak"
$ws.Range("E14").Value = "3KCO2Leak"
$ws.Range("E15").Value = "3KCO2Leak"
$ws.Range("E16").Value = "3KCO2Leak"

# AU-Mix1 / AU-Mix2 / AU-3N2O / AU-3KCO2 / AU-10N2O / AU-amb (rows 26-44) -> AU-*Leak
$ws.Range("E26").Value = "AU-Mix1Leak"
$ws.Range("E27").Value = "AU-Mix1Leak"
$ws.Range("E28").Value = "AU-Mix1Leak"
$ws.Range("E29").Value = "AU-Mix2Leak"
$ws.Range("E30").Value = "AU-Mix2Leak"
$ws.Range("E31").Value = "AU-Mix2Leak"
$ws.Range("E32").Value = "AU-3N2OLeak"
$ws.Range("E33").Value = "AU-3N2OLeak"
$ws.Range("E34").Value = "AU-3N2OLeak"
$ws.Range("E35").Value = "AU-10N2OLeak"
$ws.Range("E36").Value = "AU-10N2OLeak"
$ws.Range("E37").Value = "AU-10N2OLeak"
$ws.Range("E38").Value = "AU-3KCO2Leak"
$ws.Range("E39").Value = "AU-3KCO2Leak"
$ws.Range("E40").Value = "AU-3KCO2Leak"
$ws.Range("E41").Value = "AU-ambLeak"
$ws.Range("E42").Value = "AU-ambLeak"
$ws.Range("E43").Value = "AU-ambLeak"
$ws.Range("E44").Value = "AU-ambLeak"

# --- Move the active selection / view to where the user was working next ---
$ws.Range("E26:E44").Select()

# Try to scroll the viewport, best-effort (some hosts don't persist this to XML)
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
